$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a brand-new record as row 87 (pushes old rows 87-95 down to 88-96) ---
$ws.Rows(87).Insert()

$ws.Cells.Item(87, 1).Value = 10
$ws.Cells.Item(87, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(87, 3).Value = "La Araucanía"
$ws.Cells.Item(87, 4).Value = 44748
$ws.Cells.Item(87, 5).Value = 9
$ws.Cells.Item(87, 6).Value = 100114002
$ws.Cells.Item(87, 7).Value = "Camote"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 25
$ws.Cells.Item(87, 11).Value = 20000
$ws.Cells.Item(87, 12).Value = 20000
$ws.Cells.Item(87, 13).Value = 20000
$ws.Cells.Item(87, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(87, 15).Value = "Perú"
$ws.Cells.Item(87, 16).Value = 1000
$ws.Cells.Item(87, 17).Value = 20
$ws.Cells.Item(87, 18).Value = "Hortaliza"

# --- Insert a second brand-new record as row 94 (pushes what is currently ---
# --- rows 94-96, i.e. the originals' 93-95, down to rows 95-97) ---
$ws.Rows(94).Insert()

$ws.Cells.Item(94, 1).Value = 10
$ws.Cells.Item(94, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(94, 3).Value = "La Araucanía"
$ws.Cells.Item(94, 4).Value = 44747
$ws.Cells.Item(94, 5).Value = 9
$ws.Cells.Item(94, 6).Value = 100114002
$ws.Cells.Item(94, 7).Value = "Camote"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 20
$ws.Cells.Item(94, 11).Value = 20000
$ws.Cells.Item(94, 12).Value = 20000
$ws.Cells.Item(94, 13).Value = 20000
$ws.Cells.Item(94, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(94, 15).Value = "Perú"
$ws.Cells.Item(94, 16).Value = 1000
$ws.Cells.Item(94, 17).Value = 20
$ws.Cells.Item(94, 18).Value = "Hortaliza"
